$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update player roster table (rows 2-18) with the new data set,
# then remove the now-unused row 19 (Jaden Ivey / Detroit Pistons row).
$data = @(
  @("Fred VanVleet", "PG", "Houston Rockets"),
  @("Dennis Schröder", "PG,SG", "Golden State Warriors"),
  @("Scoot Henderson", "PG", "Portland Trail Blazers"),
  @("Amen Thompson", "SG,SF", "Houston Rockets"),
  @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
  @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
  @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
  @("Jayson Tatum", "SF,PF", "Boston Celtics"),
  @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
  @("Zach Edey", "C", "Memphis Grizzlies"),
  @("Kentavious Caldwell-Pope", "SG,SF", "Orlando Magic"),
  @("Ivica Zubac", "C", "LA Clippers"),
  @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
  @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
  @("James Harden", "PG,SG", "LA Clippers"),
  @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
  @("Jonathan Kuminga", "SF,PF", "Golden State Warriors")
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $data[$i][0]
  $ws.Cells.Item($row, 2).Value = $data[$i][1]
  $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Remove row 19, which is no longer part of the roster.
$ws.Rows.Item(19).Delete()

Write-Output "done"
